# Fix formatting of scraped floating point numbers (Importe column H)
# and normalize comma separators to periods in a handful of company names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Column H ("Importe"): these were stored as Spanish/AR-locale formatted
#        text (e.g. "1.825,00" = thousands "." + decimal ",") and must become
#        plain text with a decimal point and no thousands separator
#        (e.g. "1825.00"). H2:H199 hold one value per row, in row order.
$importeValues = @(
    "1825.00",
    "5000.00",
    "88000.00",
    "72000.00",
    "35000.00",
    "698499.80",
    "15846.80",
    "1840.00",
    "260.00",
    "55506.00",
    "144818.52",
    "181367.87",
    "16061.75",
    "1890.00",
    "4537.90",
    "512.00",
    "13355.42",
    "9810.00",
    "7000.00",
    "7677.00",
    "887.00",
    "2350.00",
    "5297.59",
    "955.96",
    "363.00",
    "38.24",
    "7596.00",
    "105.00",
    "1382.91",
    "2898.00",
    "1746.00",
    "3238.00",
    "16834.85",
    "45.00",
    "1332.80",
    "17281.68",
    "58.30",
    "5195.00",
    "281400.00",
    "5200.00",
    "930.47",
    "1950.00",
    "55001.41",
    "456.00",
    "229.20",
    "10490.00",
    "1585.33",
    "937.60",
    "5065.00",
    "591.46",
    "495.00",
    "953.00",
    "2167.00",
    "94.16",
    "354.34",
    "1100.00",
    "440.00",
    "13870.00",
    "2360.00",
    "19008.00",
    "24288.00",
    "2720.00",
    "11607.95",
    "1350.00",
    "1065.00",
    "114.72",
    "8647.00",
    "2816.35",
    "1640.00",
    "720.00",
    "2164.00",
    "300.00",
    "70.00",
    "2178.00",
    "270000.00",
    "12353.16",
    "8218.00",
    "23.15",
    "1747.00",
    "14136.15",
    "230.00",
    "199.33",
    "96.00",
    "1851.20",
    "2014.00",
    "4347.19",
    "8703.50",
    "1316.40",
    "1309.28",
    "76.04",
    "170.00",
    "12224.10",
    "1395.00",
    "4162.00",
    "1152.86",
    "1000.00",
    "138.27",
    "203.70",
    "8050.00",
    "1910.00",
    "4180.00",
    "522.00",
    "1053.00",
    "368.00",
    "130.00",
    "5500.00",
    "1585.00",
    "6000.00",
    "20000.00",
    "1500.00",
    "2946.75",
    "1393.50",
    "75.00",
    "4500.00",
    "855.54",
    "225783.00",
    "159999.87",
    "22976.16",
    "4600.00",
    "3150.00",
    "8462.00",
    "8167.50",
    "1500.00",
    "600.00",
    "1000.00",
    "6497.40",
    "5000.00",
    "950.00",
    "1750.00",
    "2250.00",
    "4655.00",
    "2450.00",
    "600.00",
    "640.00",
    "10560.00",
    "2600.00",
    "558.00",
    "1350.00",
    "1030.00",
    "380.00",
    "150.00",
    "1840.00",
    "4759.20",
    "73.00",
    "2490.00",
    "9100.00",
    "1110.00",
    "787.00",
    "4428.98",
    "27700.00",
    "2405.00",
    "180.70",
    "662.83",
    "1551.43",
    "449.86",
    "455.00",
    "6487.00",
    "7720.80",
    "4360.00",
    "272.02",
    "946.42",
    "58.78",
    "140.00",
    "2628.00",
    "110.00",
    "1655.75",
    "8217.00",
    "7374.20",
    "1273.51",
    "1449.00",
    "1778.63",
    "400.00",
    "621574.12",
    "6400.00",
    "17620.00",
    "230.00",
    "3891.51",
    "3000.00",
    "50000.00",
    "50000.00",
    "168760.92",
    "126100.00",
    "8000.00",
    "145900.00",
    "822.80",
    "1245438.22",
    "43500.00",
    "135400.00",
    "8506.90",
    "126905.00",
    "1196958.77",
    "9365.80",
    "96000.00",
    "4762.54",
    "302.01",
    "2703.00",
    "900.00",
    "1795.00"
)

$startRow = 2
for ($i = 0; $i -lt $importeValues.Length; $i++) {
    $row = $startRow + $i
    # Leading apostrophe forces Excel to keep this as literal text instead of
    # re-parsing it as a number (which would drop the decimal formatting).
    $ws.Cells.Item($row, 8).Value = "'" + $importeValues[$i]
}

# Clear the quote-prefix/text formatting picked up above so the cells keep
# their original (default) style, matching the source workbook.
$ws.Range("H2:H199").Style = "Normal"

# --- 2) A handful of "Razon social" / "Nombre Fantasia" entries used commas
#        as separators between co-contractors; replace with periods (and
#        drop the stray dots in "S.H." -> "SH") per the corrected scrape.
#        Note: in a few rows "Nombre Fantasia" (F) duplicates "Razon social"
#        (E), so both columns must be updated to stay consistent.
$ws.Range("E26").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("F26").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("E28").Value = "RAMIREZ CLAUDIA. RAMIREZ CESAR Y RAMIREZ VERONICA SH"
$ws.Range("E46").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"
$ws.Range("E69").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("F69").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("E70").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E82").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"
$ws.Range("E87").Value = "TRABICHET MARIA. VERGARA ADEL Y OTRA"
$ws.Range("F87").Value = "TRABICHET MARIA. VERGARA ADEL Y OTRA"
$ws.Range("E97").Value = "RICCOTTI. MARIANA EDITH"
$ws.Range("E146").Value = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
$ws.Range("E159").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"
$ws.Range("F109").Value = "MERCANZINI. GASTON ARIEL"
